# Sheet2 of the workbook had its "Sno." column (A) and "Acceptance Criteria"
# column (B) swapped: Acceptance Criteria moved to column A, and the Sno.
# column moved to column B and was renamed "Scenario No.".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Swap the entire columns A and B - this carries along values, number/cell
# formatting (styles) and column width, matching how the author re-ordered
# the two columns in Excel.
$ws.Columns("B").Cut()
$ws.Columns("A").Insert()

# The engine's column Cut/Insert doesn't relocate merged-cell ranges, so
# move the merges that used to group the "Acceptance Criteria" rows from
# column B back onto column A by hand.
$ws.Range("B3:B8").UnMerge()
$ws.Range("A3:A8").Merge()

$ws.Range("B9:B11").UnMerge()
$ws.Range("A9:A11").Merge()

$ws.Range("B12:B17").UnMerge()
$ws.Range("A12:A17").Merge()

$ws.Range("B18:B20").UnMerge()
$ws.Range("A18:A20").Merge()

$ws.Range("B21:B26").UnMerge()
$ws.Range("A21:A26").Merge()

$ws.Range("B27:B29").UnMerge()
$ws.Range("A27:A29").Merge()

$ws.Range("B30:B31").UnMerge()
$ws.Range("A30:A31").Merge()

$ws.Range("B33:B34").UnMerge()
$ws.Range("A33:A34").Merge()

# Column B's header used to read "Sno." -- rename it to "Scenario No."
$ws.Range("B2").Value = "Scenario No."

# Match the new view state saved with the sheet: slightly reduced zoom and
# the active selection moved to C1.
$ws.Activate()
$excel.ActiveWindow.Zoom = 86
$ws.Range("C1").Select()
